$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing existing rows 20-25 down to 21-26
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new record (same fixed fields as its neighbours,
# new Fecha/Variedad/Volumen/Origen values per the diff)
$ws.Cells.Item(20, 1).Value = 9
$ws.Cells.Item(20, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44523
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100103
$ws.Cells.Item(20, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(20, 9).Value = 100103003
$ws.Cells.Item(20, 10).Value = "Damasco"
$ws.Cells.Item(20, 11).Value = "Castle Brite"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 320
$ws.Cells.Item(20, 14).Value = 10000
$ws.Cells.Item(20, 15).Value = 10000
$ws.Cells.Item(20, 16).Value = 10000
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(20, 19).Value = 1000
$ws.Cells.Item(20, 20).Value = 10

# Match the date-style formatting used by the rest of column D
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
